# Apply the "Add files via upload" edit:
#  - rename "Mapping Respuestas" -> "Respuestas Index"
#  - add a new "Fillout Codes" sheet after it, with a question-text -> short-code
#    lookup table (13 rows), reusing the existing question-text shared strings
#    and adding 13 new short-code strings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Mapping Respuestas" -> "Respuestas Index"
# ---------------------------------------------------------------------------
$idxSheet = $wb.Worksheets.Item("Mapping Respuestas")
$idxSheet.Name = "Respuestas Index"

# ---------------------------------------------------------------------------
# 2. Add the new "Fillout Codes" sheet right after "Respuestas Index"
# ---------------------------------------------------------------------------
$codesSheet = $wb.Worksheets.Add($null, $idxSheet)
$codesSheet.Name = "Fillout Codes"

# Column widths matching the rest of the workbook's look & feel.
# (Excel's ColumnWidth property is in "characters"; the stored OOXML <col>
# width is ColumnWidth + ~5/6 character of padding, so back that padding out
# to land on the same nice round numbers used elsewhere in the workbook.)
$padding = 5.0 / 6.0
$codesSheet.Columns.Item(1).ColumnWidth = 46.0 - $padding
$codesSheet.Columns.Item(2).ColumnWidth = 16.0 - $padding

# ---------------------------------------------------------------------------
# 3. Formatting: column A reuses the same "question label" look already used
#    for column B on the "Respuestas Index" sheet; column B reuses the same
#    plain look already used for the (otherwise empty) helper cell
#    Vector Ponderacion!E6, minus its 3-decimal number format.
# ---------------------------------------------------------------------------
$vecSheet = $wb.Worksheets.Item("Vector Ponderacion")

$idxSheet.Range("B2").Copy()
$codesSheet.Range("A1:A13").PasteSpecial(-4122)

$vecSheet.Range("E6").Copy()
$codesSheet.Range("B1:B13").PasteSpecial(-4122)
$codesSheet.Range("B1:B13").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 4. Values: question text (column A, reusing existing shared strings) and
#    the new short-code identifiers (column B, new shared strings).
# ---------------------------------------------------------------------------
$rows = @(
    @("What is the size of your company?", "q1_size"),
    @("What is your annual budget?", "q2_budget"),
    @("Preferred Pricing Structure", "q3_pricing"),
    @("How soon do you need to implement the tool?", "q4_time"),
    @("What is the relationship of the tool to Excel in your use case?", "q5_excel"),
    @("What systems do you need to integrate with the FP&A tool?", "q6_integrations"),
    @("What are the key use cases for which you need the tool?", "q7_cases"),
    @("Reporting needs", "q8_reporting"),
    @("What level of modeling capability do you need?", "q9_modeling"),
    @("Collaboration and workflows needs", "q10_collaboration"),
    @("What IA capabilities are important to you?", "q11_ia"),
    @("Scalability and ability to handle large volumes of data", "q12_scalability"),
    @("Preferred implementation type", "q13_implementation")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $pair = $rows[$i]
    $codesSheet.Cells.Item($r, 1).Value = $pair[0]
    $codesSheet.Cells.Item($r, 2).Value = $pair[1]
}
